# "Generate Report for Archive"
#
# The handback/localization-status report is regenerated: the status text
# "Ready for handoff" becomes "In Translation" everywhere it appears
# (Overview sheet's per-language status columns, and the "Status" column on
# each per-language detail sheet). Because the new status text is shorter,
# the status columns that were sized to fit it shrink accordingly.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: E (zh-cn) and F (de-de) status columns, rows 2-3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- Per-language detail sheets: column C ("Status"), rows 2-3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Re-fit the status columns now that the text is shorter ---
# (stored widths are quantized to whole pixels by the host, same as real
# Excel COM; this lands on the closest achievable width to the refreshed
# auto-fit size.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
